# Update column F (dSF) values for the specified rows to reflect the
# repulled data / mean calculation described in the commit message.
# Column E (dS0) is left unchanged; only column F (dSF) is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = -2
    8  = 3
    9  = -1
    10 = -1
    11 = -4
    12 = 2
    13 = 1
    14 = -4
    15 = 3
    16 = 1
    17 = 2
    18 = -2
    22 = -7
    25 = -2
    26 = 2
    27 = -5
    28 = 2
    34 = 2
    35 = -3
    36 = -10
    38 = -4
    39 = -6
    40 = -10
    41 = -6
    43 = 3
    45 = -3
    49 = -1
    50 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
